# Update "想去人数" (number of people interested) figures on the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 89
$ws1.Range("F3").Value = 4035
$ws1.Range("F4").Value = 2363
$ws1.Range("F5").Value = 470
$ws1.Range("F8").Value = 24
$ws1.Range("F11").Value = 72
$ws1.Range("F12").Value = 131
$ws1.Range("F13").Value = 1506
$ws1.Range("F14").Value = 267
$ws1.Range("F15").Value = 2855
$ws1.Range("F16").Value = 196

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 89
$ws4.Range("F3").Value = 4035
$ws4.Range("F4").Value = 2363
$ws4.Range("F5").Value = 470
$ws4.Range("F8").Value = 24
$ws4.Range("F12").Value = 72
$ws4.Range("F13").Value = 131
$ws4.Range("F16").Value = 1506
$ws4.Range("F17").Value = 267
$ws4.Range("F18").Value = 2855
$ws4.Range("F19").Value = 196
